$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two discontinued games (no longer reported by the scraper):
# "PUMPKIN SPICE" (row 24) and "STAR PLATINUM PLAY" (row 39, which becomes
# row 38 after the first deletion shifts the remaining rows up).
$ws.Rows.Item(24).EntireRow.Delete()
$ws.Rows.Item(38).EntireRow.Delete()

# Update a couple of top-prizes-remaining counts that changed between
# scrapes.
$ws.Range("E13").Value = 4
$ws.Range("E23").Value = 62

# Refresh the "LAST SCRAPE DATE" column for every row scraped on
# 2019-03-07 to the new scrape date, 2019-03-12 (rows with other scrape
# dates, e.g. 2019-02-10 / 2019-02-19, are left untouched). Force the
# column to text first so Excel doesn't reinterpret the new value as a
# date serial number, then clear the temporary formatting so the cells
# fall back to the default style.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row
$dateRange = $ws.Range("F2:F" + $lastRow)
$dateRange.NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "2019-03-07") {
        $cell.Value = "2019-03-12"
    }
}
$dateRange.ClearFormats()
